$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (glass.csv): clear the "Num of categories" value in D2
$ws.Range("D2").ClearContents()

# Row 3 (pulsar_stars.csv): clear D3
$ws.Range("D3").ClearContents()

# Row 4 (adult.csv): update B4/C4, clear D4, add E4/F4 values
$ws.Range("B4").Value = 0.81741253051261098
$ws.Range("C4").Value = 0.45598047192839702
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = 0.812042310821806
$ws.Range("F4").Value = 0.83970707892595597

# Row 5 (heart.csv): clear D5
$ws.Range("D5").ClearContents()

# Header row: move "Clustered comonotonicity" to E1, "Weighted avg..." to F1, clear D1
$ws.Range("E1").Value = "Clustered comonotonicity"
$ws.Range("F1").Value = "Weighted avg of Naïve Bayes & Comonotonicity"
$ws.Range("D1").ClearContents()

# Column widths (engine adds a constant 5/6 padding on top of ColumnWidth when
# round-tripping through OOXML, so compensate to hit the target raw widths)
$ws.Columns.Item(4).ColumnWidth = 39.666666666666664
$ws.Columns.Item(5).ColumnWidth = 37.666666666666664
$ws.Columns.Item(6).ColumnWidth = 37.330729166666664

# Selection
$ws.Range("D4").Select()
